{"js": "// The template's final section ends with:\n//   ...<w:p> (ind left=1800 hanging=360, contains a <w:br w:type=\"page\"/> run) </w:p>\n//   <w:p pStyle=\"Bibliography\" numPr(ilvl=0,numId=0) ind(left=1080 hanging=1080)> </w:p>  <-- empty, trailing\n//   <w:sectPr> ... </w:sectPr>\n// The edit removes that trailing, content-less \"Bibliography\" paragraph so the\n// page-break paragraph above it becomes the last paragraph of the body again.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst count = paragraphs.items.length;\nif (count >= 2) {\n  const pageBreakPara = paragraphs.items[count - 2];\n  const trailingPara = paragraphs.items[count - 1];\n\n  trailingPara.load(\"text,style\");\n  await context.sync();\n\n  // Only touch it when it really is that empty trailing paragraph (defensive\n  // guard in case the document shape ever differs from what we expect).\n  if (trailingPara.text === \"\" && trailingPara.style === \"Bibliography\") {\n    // Word will not let the very last paragraph of the body (the one that\n    // owns the section break) be deleted outright, since some paragraph has\n    // to keep terminating the section. So instead we: (1) strip the trailing\n    // paragraph's own distinguishing formatting (style/numbering/indent) so\n    // it matches the plain paragraph that precedes it, then (2) delete the\n    // paragraph mark that separates the two paragraphs, which merges the\n    // page-break paragraph's content forward into the (now reformatted)\n    // trailing paragraph -- leaving a single paragraph with the page break.\n    trailingPara.style = \"Normal\";\n    trailingPara.leftIndent = pageBreakPara.leftIndent;\n    trailingPara.firstLineIndent = pageBreakPara.firstLineIndent;\n    await context.sync();\n\n    const mergeRange = pageBreakPara.getRange(\"End\").expandTo(trailingPara.getRange(\"End\"));\n    mergeRange.delete();\n    await context.sync();\n  }\n}\n", "ps1": "# The template's final section ends with:\n#   ...<w:p> (ind left=1800 hanging=360, contains a <w:br w:type=\"page\"/> run) </w:p>\n#   <w:p pStyle=\"Bibliography\" numPr(ilvl=0,numId=0) ind(left=1080 hanging=1080)> </w:p>  <-- empty, trailing\n#   <w:sectPr> ... </w:sectPr>\n# The edit removes that trailing, content-less \"Bibliography\" paragraph so the\n# page-break paragraph above it becomes the last paragraph of the body again.\n$d = $word.ActiveDocument\n\n$count = $d.Paragraphs.Count\nif ($count -ge 2) {\n    $pageBreakPara = $d.Paragraphs.Item($count - 1)\n    $trailingPara = $d.Paragraphs.Item($count)\n\n    # Paragraph.Range.Text includes the trailing paragraph-mark character\n    # (carriage return), so strip it before checking for \"no real content\".\n    $trailingText = $trailingPara.Range.Text -replace \"[\\r\\a]\", \"\"\n    $trailingStyle = $trailingPara.Style.NameLocal\n\n    # Only touch it when it really is that empty trailing paragraph (defensive\n    # guard in case the document shape ever differs from what we expect).\n    if ($trailingText -eq \"\" -and $trailingStyle -eq \"Bibliography\") {\n        # Word will not let the very last paragraph of the body (the one that\n        # owns the section break) be deleted outright, since some paragraph has\n        # to keep terminating the section. So instead we: (1) strip the trailing\n        # paragraph's own distinguishing formatting (style/numbering/indent) so\n        # it matches the plain paragraph that precedes it, then (2) delete the\n        # paragraph mark that separates the two paragraphs, which merges the\n        # page-break paragraph's content forward into the (now reformatted)\n        # trailing paragraph -- leaving a single paragraph with the page break.\n        $trailingPara.Style = \"Normal\"\n        $trailingPara.LeftIndent = $pageBreakPara.LeftIndent\n        $trailingPara.FirstLineIndent = $pageBreakPara.FirstLineIndent\n\n        $mergeRange = $pageBreakPara.Range.Duplicate\n        $mergeRange.Start = $pageBreakPara.Range.End\n        $mergeRange.End = $trailingPara.Range.End\n        $mergeRange.Delete()\n    }\n}\n"}
